# Add a new "2022" column (Q) to the table, mirroring the formatting of
# column P (the "2021" column) row by row, then fill in the new figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy each row's P-cell formatting into the corresponding Q-cell first,
# so the new column inherits the same styles as column P.
for ($r = 3; $r -le 25; $r++) {
    $ws.Range("P$r").Copy()
    $ws.Range("Q$r").PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

# Now populate the new values for column Q.
$ws.Range("Q4").Value = 2022
$ws.Range("Q5").Value = 8725
$ws.Range("Q7").Value = 8347
$ws.Range("Q8").Value = 378

for ($r = 10; $r -le 25; $r++) {
    $ws.Range("Q$r").Value = "…"
}

# Move the active selection from Q4 to Q3, matching the updated view state.
$ws.Range("Q3").Select()
